# edit.ps1 -- apply the 2025-10-05 Betfair Back/Lay odds-sheet update
#
# The source diff does two things to Sheet1:
#   1. Two match rows are dropped (the old row 10 "Uruguayan Segunda Division"
#      and old row 11 "US MLS" fixtures move up to become rows 8 and 9 with
#      refreshed odds; the trailing two physical rows disappear). We model
#      this by deleting rows 10:11 outright, which also shrinks the sheet's
#      used range from A1:AO11 down to A1:AO9 automatically.
#   2. Every remaining match row (3-9) gets a refreshed League / Time / Home /
#      Away label plus new odds across columns F:AO (a live odds snapshot).
#      Column B (Date) is identical before/after in the diff, so it is
#      intentionally left untouched -- COM auto-converts a literal
#      "2025-10-05" string into a date serial, which would NOT reproduce the
#      inlineStr cell the workbook already has there.
#
# Row 2 (Venezuelan Primera Division) is completely unchanged and so is left
# alone as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Matches for "Peruvian Segunda Division" (old row 3) and "US USL League One"
# (old row 4) are gone from the bottom of the table; rows 10 and 11 (the old
# "Uruguayan Segunda Division" / "US MLS" rows) are deleted so every other
# row shifts/refreshes into its final 3-9 position and the sheet dimension
# recomputes to A1:AO9.
$ws.Rows("10:11").Delete()

# Row 3: US USL League One -- Union Omaha vs FC Naples (18:00:00)
$ws.Cells.Item(3, 1).Value = "US USL League One"
$ws.Cells.Item(3, 3).Value = "18:00:00"
$ws.Cells.Item(3, 4).Value = "Union Omaha"
$ws.Cells.Item(3, 5).Value = "FC Naples"
$ws.Cells.Item(3, 6).Value = 1.02
$ws.Cells.Item(3, 7).Value = 680
$ws.Cells.Item(3, 8).Value = 1.01
$ws.Cells.Item(3, 9).Value = 810
$ws.Cells.Item(3, 10).Value = 1.01
$ws.Cells.Item(3, 11).Value = 680
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 1.02
$ws.Cells.Item(3, 15).Value = 1.05
$ws.Cells.Item(3, 16).Value = 1.19
$ws.Cells.Item(3, 17).Value = 1.1
$ws.Cells.Item(3, 18).Value = 1.1
$ws.Cells.Item(3, 19).Value = 1.1
$ws.Cells.Item(3, 20).Value = 1.48
$ws.Cells.Item(3, 21).Value = 1.55
$ws.Cells.Item(3, 29).Value = 980
$ws.Cells.Item(3, 30).Value = 980
$ws.Cells.Item(3, 33).Value = 990

# Row 4: Argentinian Primera Division -- Boca Juniors vs Newells (19:00:00)
$ws.Cells.Item(4, 1).Value = "Argentinian Primera Division"
$ws.Cells.Item(4, 3).Value = "19:00:00"
$ws.Cells.Item(4, 4).Value = "Boca Juniors"
$ws.Cells.Item(4, 5).Value = "Newells"
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 0
$ws.Cells.Item(4, 22).Value = 0
$ws.Cells.Item(4, 23).Value = 0
$ws.Cells.Item(4, 24).Value = 0
$ws.Cells.Item(4, 25).Value = 0
$ws.Cells.Item(4, 26).Value = 0
$ws.Cells.Item(4, 27).Value = 0
$ws.Cells.Item(4, 28).Value = 0
$ws.Cells.Item(4, 29).Value = 0
$ws.Cells.Item(4, 30).Value = 0
$ws.Cells.Item(4, 31).Value = 0
$ws.Cells.Item(4, 32).Value = 0
$ws.Cells.Item(4, 33).Value = 0
$ws.Cells.Item(4, 34).Value = 0
$ws.Cells.Item(4, 35).Value = 0
$ws.Cells.Item(4, 36).Value = 0
$ws.Cells.Item(4, 37).Value = 0
$ws.Cells.Item(4, 38).Value = 0
$ws.Cells.Item(4, 39).Value = 0
$ws.Cells.Item(4, 40).Value = 0
$ws.Cells.Item(4, 41).Value = 0

# Row 5: Argentinian Primera Division -- Rosario Central vs River Plate (21:15:00)
$ws.Cells.Item(5, 3).Value = "21:15:00"
$ws.Cells.Item(5, 4).Value = "Rosario Central"
$ws.Cells.Item(5, 5).Value = "River Plate"
$ws.Cells.Item(5, 6).Value = 1.13
$ws.Cells.Item(5, 7).Value = 1.15
$ws.Cells.Item(5, 8).Value = 90
$ws.Cells.Item(5, 9).Value = 140
$ws.Cells.Item(5, 10).Value = 8.6
$ws.Cells.Item(5, 11).Value = 9.6
$ws.Cells.Item(5, 18).Value = 1.81
$ws.Cells.Item(5, 19).Value = 2.22
$ws.Cells.Item(5, 22).Value = 1.01
$ws.Cells.Item(5, 23).Value = 8
$ws.Cells.Item(5, 24).Value = 1000
$ws.Cells.Item(5, 25).Value = 1000
$ws.Cells.Item(5, 26).Value = 1000
$ws.Cells.Item(5, 27).Value = 1000
$ws.Cells.Item(5, 28).Value = 1000
$ws.Cells.Item(5, 29).Value = 1000
$ws.Cells.Item(5, 30).Value = 1000
$ws.Cells.Item(5, 31).Value = 1000
$ws.Cells.Item(5, 32).Value = 1000
$ws.Cells.Item(5, 33).Value = 1.85
$ws.Cells.Item(5, 34).Value = 9.8
$ws.Cells.Item(5, 35).Value = 150
$ws.Cells.Item(5, 36).Value = 1000
$ws.Cells.Item(5, 37).Value = 4.6
$ws.Cells.Item(5, 38).Value = 25
$ws.Cells.Item(5, 39).Value = 350
$ws.Cells.Item(5, 40).Value = 19.5
$ws.Cells.Item(5, 41).Value = 1000

# Row 6: Colombian Primera A -- Boyaca Chico vs Atletico Nacional Medellin (21:30:00)
$ws.Cells.Item(6, 1).Value = "Colombian Primera A"
$ws.Cells.Item(6, 3).Value = "21:30:00"
$ws.Cells.Item(6, 4).Value = "Boyaca Chico"
$ws.Cells.Item(6, 5).Value = "Atletico Nacional Medellin"
$ws.Cells.Item(6, 6).Value = 2
$ws.Cells.Item(6, 7).Value = 2.02
$ws.Cells.Item(6, 8).Value = 5.9
$ws.Cells.Item(6, 9).Value = 6
$ws.Cells.Item(6, 10).Value = 2.94
$ws.Cells.Item(6, 11).Value = 3.05
$ws.Cells.Item(6, 14).Value = 2.96
$ws.Cells.Item(6, 15).Value = 1.5
$ws.Cells.Item(6, 16).Value = 1.39
$ws.Cells.Item(6, 17).Value = 3.5
$ws.Cells.Item(6, 18).Value = 1.09
$ws.Cells.Item(6, 19).Value = 11.5
$ws.Cells.Item(6, 20).Value = 1.62
$ws.Cells.Item(6, 21).Value = 2.42
$ws.Cells.Item(6, 22).Value = 1.2
$ws.Cells.Item(6, 23).Value = 1.94
$ws.Cells.Item(6, 28).Value = 2.96
$ws.Cells.Item(6, 29).Value = 3.45
$ws.Cells.Item(6, 30).Value = 8.4
$ws.Cells.Item(6, 31).Value = 48
$ws.Cells.Item(6, 32).Value = 14.5
$ws.Cells.Item(6, 33).Value = 15.5
$ws.Cells.Item(6, 34).Value = 46
$ws.Cells.Item(6, 35).Value = 220
$ws.Cells.Item(6, 36).Value = 180
$ws.Cells.Item(6, 37).Value = 190
$ws.Cells.Item(6, 38).Value = 530
$ws.Cells.Item(6, 41).Value = 240

# Row 7: Mexican Liga MX -- Pumas UNAM vs Guadalajara (22:00:00)
$ws.Cells.Item(7, 1).Value = "Mexican Liga MX"
$ws.Cells.Item(7, 3).Value = "22:00:00"
$ws.Cells.Item(7, 4).Value = "Pumas UNAM"
$ws.Cells.Item(7, 5).Value = "Guadalajara"
$ws.Cells.Item(7, 6).Value = 2.84
$ws.Cells.Item(7, 7).Value = 2.92
$ws.Cells.Item(7, 8).Value = 2.92
$ws.Cells.Item(7, 9).Value = 3.05
$ws.Cells.Item(7, 10).Value = 3.15
$ws.Cells.Item(7, 11).Value = 3.2
$ws.Cells.Item(7, 12).Value = 1.79
$ws.Cells.Item(7, 13).Value = 1.12
$ws.Cells.Item(7, 14).Value = 2.84
$ws.Cells.Item(7, 15).Value = 1.51
$ws.Cells.Item(7, 16).Value = 1.61
$ws.Cells.Item(7, 17).Value = 2.52
$ws.Cells.Item(7, 18).Value = 1.22
$ws.Cells.Item(7, 19).Value = 5.3
$ws.Cells.Item(7, 20).Value = 2.06
$ws.Cells.Item(7, 21).Value = 1.79
$ws.Cells.Item(7, 22).Value = 1.49
$ws.Cells.Item(7, 23).Value = 1.52
$ws.Cells.Item(7, 24).Value = 9.2
$ws.Cells.Item(7, 25).Value = 14.5
$ws.Cells.Item(7, 26).Value = 18
$ws.Cells.Item(7, 27).Value = 50
$ws.Cells.Item(7, 28).Value = 14
$ws.Cells.Item(7, 29).Value = 9.6
$ws.Cells.Item(7, 30).Value = 990
$ws.Cells.Item(7, 31).Value = 1000
$ws.Cells.Item(7, 33).Value = 990
$ws.Cells.Item(7, 34).Value = 990
$ws.Cells.Item(7, 35).Value = 1000
$ws.Cells.Item(7, 37).Value = 1000
$ws.Cells.Item(7, 38).Value = 1000
$ws.Cells.Item(7, 39).Value = 330
$ws.Cells.Item(7, 40).Value = 1000
$ws.Cells.Item(7, 41).Value = 1000

# Row 8: Uruguayan Segunda Division -- Central Espanol vs Atenas (22:00:00)
$ws.Cells.Item(8, 1).Value = "Uruguayan Segunda Division"
$ws.Cells.Item(8, 3).Value = "22:00:00"
$ws.Cells.Item(8, 4).Value = "Central Espanol"
$ws.Cells.Item(8, 5).Value = "Atenas"
$ws.Cells.Item(8, 6).Value = 7.4
$ws.Cells.Item(8, 7).Value = 8
$ws.Cells.Item(8, 8).Value = 1.73
$ws.Cells.Item(8, 9).Value = 1.76
$ws.Cells.Item(8, 10).Value = 3.3
$ws.Cells.Item(8, 11).Value = 3.5
$ws.Cells.Item(8, 14).Value = 3.4
$ws.Cells.Item(8, 15).Value = 1.38
$ws.Cells.Item(8, 16).Value = 1.52
$ws.Cells.Item(8, 17).Value = 2.84
$ws.Cells.Item(8, 18).Value = 1.13
$ws.Cells.Item(8, 19).Value = 8
$ws.Cells.Item(8, 20).Value = 1.74
$ws.Cells.Item(8, 21).Value = 2.16
$ws.Cells.Item(8, 22).Value = 2.3
$ws.Cells.Item(8, 23).Value = 1.14
$ws.Cells.Item(8, 25).Value = 3.7
$ws.Cells.Item(8, 26).Value = 10
$ws.Cells.Item(8, 27).Value = 55
$ws.Cells.Item(8, 28).Value = 1000
$ws.Cells.Item(8, 29).Value = 5.2
$ws.Cells.Item(8, 30).Value = 11.5
$ws.Cells.Item(8, 31).Value = 55
$ws.Cells.Item(8, 32).Value = 1000
$ws.Cells.Item(8, 33).Value = 14
$ws.Cells.Item(8, 34).Value = 30
$ws.Cells.Item(8, 35).Value = 140
$ws.Cells.Item(8, 36).Value = 1000
$ws.Cells.Item(8, 37).Value = 75
$ws.Cells.Item(8, 38).Value = 140
$ws.Cells.Item(8, 39).Value = 1000
$ws.Cells.Item(8, 40).Value = 280
$ws.Cells.Item(8, 41).Value = 140

# Row 9: US MLS -- Los Angeles FC vs Atlanta Utd (22:10:00)
$ws.Cells.Item(9, 1).Value = "US MLS"
$ws.Cells.Item(9, 3).Value = "22:10:00"
$ws.Cells.Item(9, 4).Value = "Los Angeles FC"
$ws.Cells.Item(9, 5).Value = "Atlanta Utd"
$ws.Cells.Item(9, 6).Value = 1.49
$ws.Cells.Item(9, 7).Value = 1.5
$ws.Cells.Item(9, 8).Value = 10
$ws.Cells.Item(9, 9).Value = 11
$ws.Cells.Item(9, 10).Value = 4.2
$ws.Cells.Item(9, 11).Value = 4.4
$ws.Cells.Item(9, 12).Value = 4.7
$ws.Cells.Item(9, 13).Value = 1.12
$ws.Cells.Item(9, 14).Value = 2.52
$ws.Cells.Item(9, 15).Value = 1.64
$ws.Cells.Item(9, 16).Value = 1.5
$ws.Cells.Item(9, 17).Value = 2.94
$ws.Cells.Item(9, 18).Value = 1.16
$ws.Cells.Item(9, 19).Value = 6.6
$ws.Cells.Item(9, 20).Value = 3.05
$ws.Cells.Item(9, 21).Value = 1.39
$ws.Cells.Item(9, 22).Value = 1.1
$ws.Cells.Item(9, 23).Value = 2.96
$ws.Cells.Item(9, 24).Value = 9.2
$ws.Cells.Item(9, 25).Value = 22
$ws.Cells.Item(9, 26).Value = 120
$ws.Cells.Item(9, 28).Value = 4.9
$ws.Cells.Item(9, 29).Value = 11.5
$ws.Cells.Item(9, 30).Value = 60
$ws.Cells.Item(9, 31).Value = 450
$ws.Cells.Item(9, 32).Value = 6.4
$ws.Cells.Item(9, 33).Value = 13
$ws.Cells.Item(9, 34).Value = 75
$ws.Cells.Item(9, 36).Value = 13
$ws.Cells.Item(9, 37).Value = 28
$ws.Cells.Item(9, 38).Value = 160
$ws.Cells.Item(9, 40).Value = 15
$ws.Cells.Item(9, 41).Value = 1000
